# The workbook contains a time series table (years in column A, coefficients
# in columns B:Y) for rows 2-5 (2007年, 2010年, 2012年, 2015年).
# This edit removes the 2007年 row (row 2), shifting the remaining rows
# (2010年, 2012年, 2015年) up by one, and shrinking the used range from
# A1:Y5 to A1:Y4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
